$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.590.71"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.839.16"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("D4").Value = "'0.9981"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'319.17"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.5335"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("D8").Value = "'0.4013"
$ws.Range("E8").Value = "  +6.17%  "
$ws.Range("D9").Value = "'0.07602"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "'41.88"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'1.114"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'6.324"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "'7.620"
$ws.Range("E13").Value = "  +4.51%  "
$ws.Range("D14").Value = "'0.9977"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "'20.87"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "1.833.76"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "'90.01"
$ws.Range("D18").Value = "'0.00001075"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "'0.06597"
$ws.Range("D20").Value = "'17.72"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'0.9986"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'6.084"
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("D23").Value = "28.585.27"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'11.25"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").Value = "'2.476"
$ws.Range("E26").Value = "  +6.84%  "
$ws.Range("D29").Value = "2.042.52"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").Value = "'124.09"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Value = "'1.130"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("D32").Value = "'0.1100"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("D33").Value = "'5.721"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'0.07247"
$ws.Range("E35").Value = "  +12.16%  "
$ws.Range("D36").Value = "'0.2266"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'5.279"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").Value = "'0.02354"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").Value = "'8.847"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").Value = "'11.40"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "'0.6306"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "'1.202"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'1.412"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").Value = "'0.9979"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'13.59"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "'0.5861"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'126.15"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "'1.995"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").Value = "'1.196"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'0.06934"
$ws.Range("E51").Value = "  +0.84%  "

# Row 27/28: swap Monero and EthereumClassic entries
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.69"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'156.94"
$ws.Range("E28").Value = "  -1.57%  "
